$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E15").Value = 143
$ws.Range("E17").Value = 95

$ws.Range("E19").Value = 45
$ws.Range("F19").Value = 23
$ws.Range("H19").Value = 23

$ws.Range("E28").Value = 12
$ws.Range("E29").Value = 16

$ws.Range("E33").Value = 34
$ws.Range("F33").Value = 11
$ws.Range("H33").Value = 11

$ws.Range("E36").Value = 86

$ws.Range("E37").Value = 44
$ws.Range("F37").Value = 22
$ws.Range("H37").Value = 22

$ws.Range("E41").Value = 32

$ws.Range("E63").Value = 28
$ws.Range("F63").Value = 9
$ws.Range("H63").Value = 9

$ws.Range("E64").Value = 29

$ws.Range("E67").Value = 36
$ws.Range("F67").Value = 20
$ws.Range("H67").Value = 20

$ws.Range("E71").Value = 26

$ws.Range("E73").Value = 27

$ws.Range("E76").Value = 44

$ws.Range("E77").Value = 48
$ws.Range("F77").Value = 16
$ws.Range("H77").Value = 16

$ws.Range("E82").Value = 13

$wb.Save()
